# ---------------------------------------------------------------------------
# Edit summary (per the source diff):
#   1. The single table on slide 16 switches its table style (tableStyleId)
#      from {2D0B1C87-61E2-46CE-82BE-78B8B920CD8D} to
#      {0CD77837-BD8D-438D-BB2F-3ED9452FB872}.
#   2. The deck's theme (ppt/theme/theme1.xml, used by the slide master /
#      all slides) switches its colour scheme from the custom "Integral"
#      palette to the stock Office palette (dk1/lt1 unchanged; dk2, lt2,
#      accent1-6, hlink and folHlink all change). Font scheme / format
#      scheme are already identical between the two themes in this deck,
#      so only the colours need touching.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -----------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $shp = $s.Shapes.Item($shi)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{0CD77837-BD8D-438D-BB2F-3ED9452FB872}")
        }
    }
}

# --- 2. Theme colour scheme ----------------------------------------------
function Get-BGRInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Order matches MsoThemeColorSchemeIndex 1..12:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$slide1 = $p.Slides.Item(1)
$tcs = $slide1.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = Get-BGRInt $officeColors[$i - 1]
}
